$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9488468170166016
$ws.Range("B1").Value = 2.070756196975708
$ws.Range("C1").Value = 8.14568042755127
$ws.Range("D1").Value = 2.202033519744873
$ws.Range("E1").Value = 1.006226897239685
